$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Modify existing row 11: change date and values
$ws.Range("D11").Value = 44449
$ws.Range("J11").Value = 25
$ws.Range("K11").Value = 80000
$ws.Range("L11").Value = 80000
$ws.Range("M11").Value = 80000
$ws.Range("P11").Value = 3200

# Add new row 12
$ws.Range("A12").Value = 12
$ws.Range("B12").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value = 44449
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 100112021
$ws.Range("G12").Value = "Ají"
$ws.Range("H12").Value = "Americana (o)"
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 20
$ws.Range("K12").Value = 75000
$ws.Range("L12").Value = 75000
$ws.Range("M12").Value = 75000
$ws.Range("N12").Value = "$/caja 15 kilos"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 5000
$ws.Range("Q12").Value = 15
$ws.Range("R12").Value = "Hortaliza"

# Add new row 13 (restores original row 11 data)
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 44319
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100112021
$ws.Range("G13").Value = "Ají"
$ws.Range("H13").Value = "Americana (o)"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 20
$ws.Range("K13").Value = 30000
$ws.Range("L13").Value = 30000
$ws.Range("M13").Value = 30000
$ws.Range("N13").Value = "$/caja 25 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 1200
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"

# Copy date style (s="2") from D11 to D12 and D13
$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D13").PasteSpecial(-4122)
